# Coding Ninja - Finding Square root of a number
# Adds a new row (row 15) describing a "Finding Sqrt of a number using
# Binary Search" question, clears a stray value in B9, and moves the
# sheet selection onto the newly added comment cell (F15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: the "LC/CN" tag in B9 was removed (cell kept, value cleared) ---
$ws.Range("B9").ClearContents()

# --- New row 15 ----------------------------------------------------------
# Reuse the formatting already present elsewhere in the sheet for the
# matching columns, so no spurious new cell styles get created.
$ws.Range("D10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "CN/GFG"

$ws.Range("C2").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "Finding Sqrt of a number using Binary Search"

$ws.Range("D10").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "Java"

$ws.Range("E8").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "Medium"

# F15 needs a brand-new style (bordered + wrap text, no fill) together
# with a rich-text comment describing the naive & binary-search approach.
$full = "1. Naïve Approach :-`nWe will iterate till the number in for loop and check if square of i is less than N.`n2. Binary Search :-`n"
$ws.Range("F15").Value = $full
$ws.Range("F15").Borders.LineStyle = 1
$ws.Range("F15").WrapText = $true

$rng = $ws.Range("F15")
$rng.Characters(1, 20).Font.Bold = $true
$rng.Characters(21, 86).Font.Bold = $false
$rng.Characters(107, 19).Font.Bold = $true
$rng.Characters(126, 1).Font.Bold = $false

$ws.Rows.Item(15).RowHeight = 75

# --- Selection moves to the newly filled-in comment cell -----------------
$ws.Range("F15").Select() | Out-Null
